# Rename the sheet to match the regenerated export ("Sheet 1" -> "Sheet1")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Sheet1"

# Apply the header-row formatting used by the export (bold, centered)
$header = $ws.Range("A1:F1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
